$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of row 182 into the two new rows (183, 184) so the
# new rows pick up the exact same cell styles used by the existing
# "housing model" parameter block.
$ws.Range("A182:J182").Copy() | Out-Null
$ws.Range("A183:J183").PasteSpecial(-4122) | Out-Null
$ws.Range("A182:J182").Copy() | Out-Null
$ws.Range("A184:J184").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the values column-by-column (B then C for both rows, then G for
# both rows) so that any newly introduced shared strings land in the same
# order as the authored workbook.
$ws.Cells.Item(183, 2).Value = "empty_coop"
$ws.Cells.Item(183, 3).Value = "similar to wohn.modell.anteil.leerwhg"
$ws.Cells.Item(184, 2).Value = "empty_private"
$ws.Cells.Item(184, 3).Value = "similar to wohn.modell.anteil.leerwhg"
$ws.Cells.Item(183, 7).Value = "Percentage of empty apartments (cooperative housing)"
$ws.Cells.Item(184, 7).Value = "Percentage of empty apartments (private housing)"

# Remaining cells.
$ws.Cells.Item(183, 1).Value = "housing model"
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = "percent"
$ws.Cells.Item(183, 6).Value = "low"
$ws.Cells.Item(183, 8).Value = 0
$ws.Cells.Item(183, 9).Value = 0
$ws.Cells.Item(183, 10).Value = 0

$ws.Cells.Item(184, 1).Value = "housing model"
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = "percent"
$ws.Cells.Item(184, 6).Value = "low"
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 9).Value = 0
$ws.Cells.Item(184, 10).Value = 0

# Match the selection left behind in the authored workbook.
$ws.Range("A182").Select() | Out-Null
